# Apply cryptos list update (commit: "Updated cryptos list on Sun Jan  7 07:09:24 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.461.34"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.250.03"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'308.57"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").Value = "'94.82"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E10").Value = "  +4.61%  "
$ws.Range("D11").Value = "'0.0811"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").Value = "'7.27"
$ws.Range("E12").Value = "  +3.66%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.319.15"
$ws.Range("E14").Value = "  +3.12%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.840"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").Value = "'13.67"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").Value = "44.206.43"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("E18").Value = "  +2.82%  "
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "'12.27"
$ws.Range("E19").Value = "  +4.00%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.43"
$ws.Range("E20").Value = "  +6.15%  "
$ws.Range("D21").Value = "'66.12"
$ws.Range("E21").Value = "  +3.77%  "
$ws.Range("D22").Value = "'3.01"
$ws.Range("E22").Value = "  +5.20%  "
$ws.Range("D23").Value = "'237.26"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  +5.82%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +5.23%  "
$ws.Range("D27").Value = "'38.17"
$ws.Range("E27").Value = "  +7.53%  "
$ws.Range("D28").Value = "'9.87"
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").Value = "'20.10"
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("D31").Value = "'153.22"
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").Value = "'0.0799"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "'3.13"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").Value = "'0.121"
$ws.Range("E35").Value = "  +3.65%  "
$ws.Range("E36").Value = "  +5.01%  "
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("D38").Value = "'3.46"
$ws.Range("E38").Value = "  +7.99%  "
$ws.Range("D39").Value = "'14.63"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").Value = "'3.84"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").Value = "'0.0304"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "1.751.44"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("E44").Value = "  +6.67%  "
$ws.Range("D45").Value = "'80.95"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "'71.25"
$ws.Range("E46").Value = "  +6.65%  "
$ws.Range("D47").Value = "'100.00"
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("D48").Value = "'4.90"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").Value = "'1.60"
$ws.Range("E49").Value = "  +8.85%  "
$ws.Range("D50").Value = "'55.62"
$ws.Range("E50").Value = "  +5.84%  "
$ws.Range("D51").Value = "'8.19"
$ws.Range("E51").Value = "  +2.96%  "
